# Auto-generated Excel COM-interop script applying the Bahamut_Profits.xlsx diff.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H, I, J, K, L, M, N) for specific leve rows across the ALC, ARM, BSM, CUL and
# GSM sheets, matching the upstream market-price refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 50002810
$ws.Range("I76").Value = 50002810
$ws.Range("K76").Value = 50002810
$ws.Range("M76").Value = -50002495

$ws.Range("H79").Value = 50002810
$ws.Range("I79").Value = 50002810
$ws.Range("K79").Value = 50002810
$ws.Range("M79").Value = -50001718

$ws.Range("H86").Value = 55558336
$ws.Range("I86").Value = 2004.4667
$ws.Range("J86").Value = 333340000
$ws.Range("K86").Value = 2004.4667
$ws.Range("L86").Value = 333340000
$ws.Range("M86").Value = -881.4666999999999
$ws.Range("N86").Value = -333342246

$ws.Range("H89").Value = 55558336
$ws.Range("I89").Value = 2004.4667
$ws.Range("J89").Value = 333340000
$ws.Range("K89").Value = 10022.3335
$ws.Range("L89").Value = 1666700000
$ws.Range("M89").Value = -4406.333500000001
$ws.Range("N89").Value = -1666711232

$ws.Range("H116").Value = 3797.4119
$ws.Range("I116").Value = 4055
$ws.Range("J116").Value = 3429.4285
$ws.Range("K116").Value = 4055
$ws.Range("L116").Value = 3429.4285
$ws.Range("M116").Value = -613
$ws.Range("N116").Value = -10313.4285

$ws.Range("H141").Value = 7499.5
$ws.Range("I141").Value = 7499.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 22498.5
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -17318.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6471.8164
$ws.Range("I32").Value = 4163.9766
$ws.Range("K32").Value = 4163.9766
$ws.Range("M32").Value = -3876.9766

$ws.Range("H63").Value = 2155.2727
$ws.Range("I63").Value = 2126.5806
$ws.Range("K63").Value = 2126.5806
$ws.Range("M63").Value = -1440.5806

$ws.Range("H66").Value = 2155.2727
$ws.Range("I66").Value = 2126.5806
$ws.Range("K66").Value = 10632.903
$ws.Range("M66").Value = -7200.902999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5068.75
$ws.Range("I105").Value = 4282.65
$ws.Range("K105").Value = 4282.65
$ws.Range("M105").Value = -2535.65

$ws.Range("H107").Value = 13762
$ws.Range("I107").Value = 1509.1111
$ws.Range("K107").Value = 1509.1111
$ws.Range("M107").Value = 410.8888999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1777.7778
$ws.Range("J49").Value = 3333.3333
$ws.Range("L49").Value = 9999.999899999999
$ws.Range("N49").Value = -10311.9999

$ws.Range("H58").Value = 3712.9167
$ws.Range("I58").Value = 1647.5
$ws.Range("J58").Value = 4126
$ws.Range("K58").Value = 4942.5
$ws.Range("L58").Value = 12378
$ws.Range("M58").Value = -4814.5
$ws.Range("N58").Value = -12634

$ws.Range("H63").Value = 4238.875
$ws.Range("I63").Value = 1011
$ws.Range("J63").Value = 4700
$ws.Range("K63").Value = 3033
$ws.Range("L63").Value = 14100
$ws.Range("M63").Value = -2284
$ws.Range("N63").Value = -15598

$ws.Range("H64").Value = 1864.5454
$ws.Range("J64").Value = 1931
$ws.Range("L64").Value = 5793
$ws.Range("N64").Value = -6333

$ws.Range("H66").Value = 4238.875
$ws.Range("I66").Value = 1011
$ws.Range("J66").Value = 4700
$ws.Range("K66").Value = 9099
$ws.Range("L66").Value = 42300
$ws.Range("M66").Value = -5355
$ws.Range("N66").Value = -49788

$ws.Range("H67").Value = 1864.5454
$ws.Range("J67").Value = 1931
$ws.Range("L67").Value = 5793
$ws.Range("N67").Value = -7665

$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 6000
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -5617
$ws.Range("N76").Value = -15766

$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -4674
$ws.Range("N79").Value = -17652

$ws.Range("H81").Value = 2500
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H82").Value = 12394.25
$ws.Range("J82").Value = 12394.25
$ws.Range("L82").Value = 37182.75
$ws.Range("N82").Value = -37994.75

$ws.Range("H84").Value = 2500
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H85").Value = 12394.25
$ws.Range("J85").Value = 12394.25
$ws.Range("L85").Value = 37182.75
$ws.Range("N85").Value = -39990.75

$ws.Range("H86").Value = 516.6667
$ws.Range("I86").Value = 400
$ws.Range("J86").Value = 550
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 1650
$ws.Range("M86").Value = -14
$ws.Range("N86").Value = -4022

$ws.Range("H87").Value = 14054.444
$ws.Range("I87").Value = 3000
$ws.Range("J87").Value = 17212.857
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 51638.571
$ws.Range("M87").Value = -7752
$ws.Range("N87").Value = -54134.571

$ws.Range("H88").Value = 8000
$ws.Range("J88").Value = 8000
$ws.Range("L88").Value = 24000
$ws.Range("N88").Value = -24856

$ws.Range("H89").Value = 516.6667
$ws.Range("I89").Value = 400
$ws.Range("J89").Value = 550
$ws.Range("K89").Value = 3600
$ws.Range("L89").Value = 4950
$ws.Range("M89").Value = 2328
$ws.Range("N89").Value = -16806

$ws.Range("H90").Value = 14054.444
$ws.Range("I90").Value = 3000
$ws.Range("J90").Value = 17212.857
$ws.Range("K90").Value = 27000
$ws.Range("L90").Value = 154915.713
$ws.Range("M90").Value = -20760
$ws.Range("N90").Value = -167395.713

$ws.Range("H91").Value = 8000
$ws.Range("J91").Value = 8000
$ws.Range("L91").Value = 24000
$ws.Range("N91").Value = -26964

$ws.Range("H110").Value = 8674.4375
$ws.Range("I110").Value = 3465.1667
$ws.Range("K110").Value = 10395.5001
$ws.Range("M110").Value = -6305.500100000001

$ws.Range("H122").Value = 652.24445
$ws.Range("I122").Value = 539.8
$ws.Range("J122").Value = 666.3
$ws.Range("K122").Value = 4858.2
$ws.Range("L122").Value = 5996.7
$ws.Range("M122").Value = -2408.2
$ws.Range("N122").Value = -10896.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5366.846
$ws.Range("I70").Value = 5312.5
$ws.Range("J70").Value = 5453.8
$ws.Range("K70").Value = 5312.5
$ws.Range("L70").Value = 5453.8
$ws.Range("M70").Value = -5042.5
$ws.Range("N70").Value = -5993.8

$ws.Range("H73").Value = 5366.846
$ws.Range("I73").Value = 5312.5
$ws.Range("J73").Value = 5453.8
$ws.Range("K73").Value = 5312.5
$ws.Range("L73").Value = 5453.8
$ws.Range("M73").Value = -4376.5
$ws.Range("N73").Value = -7325.8

$ws.Range("H132").Value = 3039.3572
$ws.Range("I132").Value = 2823.75
$ws.Range("K132").Value = 8471.25
$ws.Range("M132").Value = -5941.25
